$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "EmpoyeeSearch" test case, rows 9-11 (mirrors the existing
# TestCase blocks above: a bold TestCase name cell in column A,
# then Keyword/Object/ObjectType/value rows).

# Row 9: TestCase = EmpoyeeSearch, Keyword = CLICK on admin(id)
$ws.Range("A9").Value = "EmpoyeeSearch"
$ws.Range("B9").Value = "CLICK"
$ws.Range("C9").Value = "admin"
$ws.Range("D9").Value = "id"

# Row 10: SETTEXT employeename(id) = "a "
$ws.Range("B10").Value = "SETTEXT"
$ws.Range("C10").Value = "employeename"
$ws.Range("D10").Value = "id"

# Row 11: CLICK search(id)
$ws.Range("B11").Value = "CLICK"
$ws.Range("C11").Value = "search"
$ws.Range("D11").Value = "id"

# Written last so the shared-string table grows in the same order as the
# source workbook (EmpoyeeSearch, employeename, search, "a ").
$ws.Range("E10").Value = "a "

# Match the bordered-table look used by the rest of the sheet.
$ws.Range("A9:E11").Borders.LineStyle = 1

$ws.Range("E11").Select()
